$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select column G (the "Type" column) and delete it entirely, shifting
# columns to its right (Folio No, Instrument, Currency, Investment Domicile)
# one position to the left.
$col = $ws.Range("G1:G1048576")
$col.Select()
$col.EntireColumn.Delete()
